$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sessions")

# Insert two new rows above what is currently row 39 ("Batch Settings" blank
# spacer row), shifting the "Batch Settings" section (and everything below
# it) down by two rows.
$ws.Rows("39:40").Insert()

# Copy the cell formatting from row 38 ("New Vehicle Price Sales Response
# Elasticity") - a "Session Settings" value row - onto the two freshly
# inserted rows so they pick up the same styles (row style, borders, fill,
# alignment) instead of Excel's auto-generated insert formatting.
$ws.Range("A38:S38").Copy()
$ws.Range("A39:S40").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the two new "Consumer Pricing Multiplier" parameter rows.
$ws.Range("A39").Value = "Consumer Pricing Multiplier Min"
$ws.Range("B39").Value = "#"
$ws.Range("C39").Value = 0.95

$ws.Range("A40").Value = "Consumer Pricing Multiplier Max"
$ws.Range("B40").Value = "#"
$ws.Range("C40").Value = 1.05

# Update the current selection to match the saved view state (frozen panes
# stay split at row 10, same as before the edit).
$ws.Range("A39:XFD40").Select()
